$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from A44 onto A45:A46 (keeps reusing the same
# cellXfs entry instead of Excel fabricating a new custom numFmt).
$ws.Range("A44").Copy()
$ws.Range("A45:A46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 45
$ws.Range("A45").Value = 43716
$ws.Range("B45").Value = 2256.8387198168102
$ws.Range("C45").Value = 2207.0300000000002
$ws.Range("D45").Formula = "=100*(B45-C45)/C45"
$ws.Range("E45").Value = 169
$ws.Range("F45").Value = "Opened CRM (8/30/2019)"

# Row 46
$ws.Range("A46").Value = 43716
$ws.Range("B46").Value = 2221.5640712207201
$ws.Range("C46").Value = 2207.0300000000002
$ws.Range("D46").Formula = "=100*(B46-C46)/C46"
$ws.Range("E46").Value = 169
$ws.Range("F46").Value = "Opened CRM (9/8/2019)"

$ws.Range("E47").Select()
